# Add a new "2023" column (L) to the corruption-perception index table by
# copying the formatting of the existing "2022" column (K) and filling in
# the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles/number formats/etc.) from column K (rows 4-14)
# into the new column L, exactly mirroring how Excel extends a table when a
# new year column is added.
$ws.Range("K4:K14").Copy()
$ws.Range("L4:L14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New column header (year)
$ws.Cells.Item(4, 12).Value = 2023

# New column data values (2023 figures)
$ws.Cells.Item(5, 12).Value = 22.743990309495757
$ws.Cells.Item(6, 12).Value = 52.401334422687093
$ws.Cells.Item(7, 12).Value = 40.084286291781751
$ws.Cells.Item(8, 12).Value = 58.6564425462321
$ws.Cells.Item(9, 12).Value = 52.689880705632987
$ws.Cells.Item(10, 12).Value = 19.88866894869804
$ws.Cells.Item(11, 12).Value = 35.972443863264772
$ws.Cells.Item(12, 12).Value = 12.061786277026036
$ws.Cells.Item(13, 12).Value = -0.064288010286095529
$ws.Cells.Item(14, 12).Value = 34.132731805770057

# Row height adjustments (Excel auto-resized these rows when the new
# column's content was added)
$ws.Rows(1).RowHeight = 67.5
$ws.Rows("4:14").RowHeight = 14.25

# Clear the stale selection left over from editing (M7), matching the
# cleaned-up sheetView in the saved file.
$ws.Range("A1").Select()
